$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.Value = "'66.139.37"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.13%  "

# Row 3
$cell = $ws.Range("D3")
$cell.Value = "'3.158.67"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -1.58%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$cell = $ws.Range("D5")
$cell.Value = "'600.48"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.49%  "

# Row 6
$cell = $ws.Range("D6")
$cell.Value = "'153.27"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.53%  "

# Row 8
$cell = $ws.Range("D8")
$cell.Value = "'0.548"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +2.52%  "

# Row 9
$cell = $ws.Range("D9")
$cell.Value = "'3.150.74"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -1.80%  "

# Row 10
$ws.Range("E10").Value = "  -2.56%  "

# Row 11
$ws.Range("E11").Value = "  -9.63%  "

# Row 12
$cell = $ws.Range("D12")
$cell.Value = "'0.507"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.27%  "

# Row 13
$cell = $ws.Range("D13")
$cell.Value = "'0.0000264"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -3.43%  "

# Row 14
$cell = $ws.Range("D14")
$cell.Value = "'38.23"
$cell.Style = "Normal"

# Row 15
$cell = $ws.Range("D15")
$cell.Value = "'3.678.22"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -1.56%  "

# Row 16
$cell = $ws.Range("D16")
$cell.Value = "'66.237.74"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.03%  "

# Row 17
$cell = $ws.Range("D17")
$cell.Value = "'7.34"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -1.47%  "

# Row 18
$cell = $ws.Range("D18")
$cell.Value = "'3.162.73"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -1.58%  "

# Row 19
$ws.Range("E19").Value = "  +0.16%  "

# Row 20
$cell = $ws.Range("D20")
$cell.Value = "'508.15"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.59%  "

# Row 21
$cell = $ws.Range("D21")
$cell.Value = "'15.27"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -2.09%  "

# Row 22
$cell = $ws.Range("D22")
$cell.Value = "'0.724"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -1.49%  "

# Row 23
$ws.Range("E23").Value = "  +0.53%  "

# Row 24
$cell = $ws.Range("D24")
$cell.Value = "'14.51"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -4.61%  "

# Row 25
$cell = $ws.Range("D25")
$cell.Value = "'84.48"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.81%  "

# Row 27
$cell = $ws.Range("D27")
$cell.Value = "'2.98"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -1.45%  "

# Row 28
$cell = $ws.Range("D28")
$cell.Value = "'9.04"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -2.22%  "

# Row 29
$ws.Range("E29").Value = "  +5.73%  "

# Row 30
$cell = $ws.Range("D30")
$cell.Value = "'3.03"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +5.24%  "

# Row 31
$cell = $ws.Range("D31")
$cell.Value = "'6.87"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.03%  "

# Row 32
$cell = $ws.Range("D32")
$cell.Value = "'27.80"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -1.31%  "

# Row 33
$ws.Range("E33").Value = "  +0.03%  "

# Row 34
$ws.Range("E34").Value = "  -1.75%  "

# Row 35
$cell = $ws.Range("D35")
$cell.Value = "'6.46"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -2.58%  "

# Row 36
$cell = $ws.Range("D36")
$cell.Value = "'504.74"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +4.15%  "

# Row 37
$cell = $ws.Range("D37")
$cell.Value = "'54.71"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -1.15%  "

# Row 38
$cell = $ws.Range("D38")
$cell.Value = "'0.0881"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -3.51%  "

# Row 39
$cell = $ws.Range("D39")
$cell.Value = "'0.0416"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -1.28%  "

# Row 40
$ws.Range("E40").Value = "  +6.42%  "

# Row 41
$cell = $ws.Range("D41")
$cell.Value = "'8.77"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -0.73%  "

# Row 42
$cell = $ws.Range("D42")
$cell.Value = "'0.0₃0670"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +4.13%  "

# Row 43
$cell = $ws.Range("D43")
$cell.Value = "'0.294"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -1.06%  "

# Row 44
$cell = $ws.Range("D44")
$cell.Value = "'2.76"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -7.43%  "

# Row 45
$ws.Range("E45").Value = "  -4.09%  "

# Row 46
$cell = $ws.Range("D46")
$cell.Value = "'2.815.97"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -4.86%  "

# Row 47
$cell = $ws.Range("D47")
$cell.Value = "'27.75"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -3.80%  "

# Row 48
$ws.Range("E48").Value = "  -0.11%  "

# Row 49
$ws.Range("E49").Value = "  +0.40%  "

# Row 50
$ws.Range("E50").Value = "  +0.13%  "

# Row 51
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$cell = $ws.Range("D51")
$cell.Value = "'2.58"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +6.42%  "
